$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Enterprises density (per 1000 people) row: Micro, SMEs, MSMEs
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "30.45"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "0.91"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.36"

# Enterprises (% of total) row: Micro, SMEs, MSMEs
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "97.02"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "2.89"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "99.91"
